$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Report is ready for handoff: status text + refreshed timestamps ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 09:27:49"
$wsDeDe.Range("H2").Value = "2016-08-30 09:27:49"

$wsZhCn.Range("H2").Value = "2016-08-30 09:27:45"

# --- Column width adjustments: narrower date/status columns ---
# Target stored width is 17.2159881591797 (121px @ MDW=7); this runtime's
# ColumnWidth setter quantizes the stored width to steps of 1/6, so we pick
# the input that lands on the closest reachable value (17.1666...).
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
